$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H17").Value = 1249.3529
$ws_ALC.Range("I17").Value = 576
$ws_ALC.Range("J17").Value = 1365.4482
$ws_ALC.Range("K17").Value = 1728
$ws_ALC.Range("L17").Value = 4096.3446
$ws_ALC.Range("M17").Value = -1560
$ws_ALC.Range("N17").Value = -4432.3446
$ws_ALC.Range("H42").Value = 103.8
$ws_ALC.Range("I42").Value = 103.22222
$ws_ALC.Range("J42").Value = 109
$ws_ALC.Range("K42").Value = 309.66666
$ws_ALC.Range("L42").Value = 327
$ws_ALC.Range("M42").Value = -79.66665999999998
$ws_ALC.Range("N42").Value = -787
$ws_ALC.Range("H58").Value = 75008504
$ws_ALC.Range("I58").Value = 83333384
$ws_ALC.Range("J58").Value = 71440700
$ws_ALC.Range("K58").Value = 250000152
$ws_ALC.Range("L58").Value = 214322100
$ws_ALC.Range("M58").Value = -250000002
$ws_ALC.Range("N58").Value = -214322400
$ws_ALC.Range("H106").Value = 166668640
$ws_ALC.Range("I106").Value = 166668640
$ws_ALC.Range("K106").Value = 166668640
$ws_ALC.Range("M106").Value = -166668009
$ws_ALC.Range("H132").Value = 2113.3794
$ws_ALC.Range("I132").Value = 2014.3704
$ws_ALC.Range("K132").Value = 6043.1112
$ws_ALC.Range("M132").Value = -3513.1112
$ws_ALC.Range("H138").Value = 6259261
$ws_ALC.Range("I138").Value = 1523.375
$ws_ALC.Range("J138").Value = 12516999
$ws_ALC.Range("K138").Value = 4570.125
$ws_ALC.Range("L138").Value = 37550997
$ws_ALC.Range("M138").Value = 569.875
$ws_ALC.Range("N138").Value = -37561277
$ws_ALC.Range("H141").Value = 1621.1111
$ws_ALC.Range("I141").Value = 1199.25
$ws_ALC.Range("J141").Value = 4996
$ws_ALC.Range("K141").Value = 3597.75
$ws_ALC.Range("L141").Value = 14988
$ws_ALC.Range("M141").Value = 1582.25
$ws_ALC.Range("N141").Value = -25348
$ws_ARM.Range("H32").Value = 1496687.9
$ws_ARM.Range("I32").Value = 1651471.2
$ws_ARM.Range("J32").Value = 26245
$ws_ARM.Range("K32").Value = 1651471.2
$ws_ARM.Range("L32").Value = 26245
$ws_ARM.Range("M32").Value = -1651184.2
$ws_ARM.Range("N32").Value = -26819
$ws_ARM.Range("H61").Value = 5354.4595
$ws_ARM.Range("I61").Value = 1324.48
$ws_ARM.Range("K61").Value = 1324.48
$ws_ARM.Range("M61").Value = -1112.48
$ws_ARM.Range("H74").Value = 2327.4517
$ws_ARM.Range("I74").Value = 1356.1428
$ws_ARM.Range("J74").Value = 4367.2
$ws_ARM.Range("K74").Value = 1356.1428
$ws_ARM.Range("L74").Value = 4367.2
$ws_ARM.Range("M74").Value = -482.1428000000001
$ws_ARM.Range("N74").Value = -6115.2
$ws_ARM.Range("H77").Value = 2327.4517
$ws_ARM.Range("I77").Value = 1356.1428
$ws_ARM.Range("J77").Value = 4367.2
$ws_ARM.Range("K77").Value = 6780.714
$ws_ARM.Range("L77").Value = 21836
$ws_ARM.Range("M77").Value = -2412.714
$ws_ARM.Range("N77").Value = -30572
$ws_ARM.Range("H97").Value = 2874221.8
$ws_ARM.Range("J97").Value = 16667926
$ws_ARM.Range("L97").Value = 16667926
$ws_ARM.Range("N97").Value = -16668918
$ws_ARM.Range("H132").Value = 5785.636
$ws_ARM.Range("I132").Value = 3072.6843
$ws_ARM.Range("K132").Value = 9218.052899999999
$ws_ARM.Range("M132").Value = -6688.052899999999
$ws_ARM.Range("H136").Value = 5354.4595
$ws_ARM.Range("I136").Value = 1324.48
$ws_ARM.Range("K136").Value = 3973.44
$ws_ARM.Range("M136").Value = -1423.44
$ws_BSM.Range("H20").Value = 5377804.5
$ws_BSM.Range("I20").Value = 7937928
$ws_BSM.Range("J20").Value = 1545.4
$ws_BSM.Range("K20").Value = 7937928
$ws_BSM.Range("L20").Value = 1545.4
$ws_BSM.Range("M20").Value = -7937681
$ws_BSM.Range("N20").Value = -2039.4
$ws_BSM.Range("H86").Value = 55613380
$ws_BSM.Range("I86").Value = 78389.766
$ws_BSM.Range("J86").Value = 200004370
$ws_BSM.Range("K86").Value = 78389.766
$ws_BSM.Range("L86").Value = 200004370
$ws_BSM.Range("M86").Value = -77266.766
$ws_BSM.Range("N86").Value = -200006616
$ws_BSM.Range("H89").Value = 55613380
$ws_BSM.Range("I89").Value = 78389.766
$ws_BSM.Range("J89").Value = 200004370
$ws_BSM.Range("K89").Value = 391948.83
$ws_BSM.Range("L89").Value = 1000021850
$ws_BSM.Range("M89").Value = -386332.83
$ws_BSM.Range("N89").Value = -1000033082
$ws_BSM.Range("H105").Value = 2480.7
$ws_BSM.Range("I105").Value = 2256.3333
$ws_BSM.Range("K105").Value = 2256.3333
$ws_BSM.Range("M105").Value = -509.3332999999998
$ws_BSM.Range("H134").Value = 3946.192
$ws_BSM.Range("I134").Value = 1436.32
$ws_BSM.Range("K134").Value = 4308.96
$ws_BSM.Range("M134").Value = -1773.96
$ws_CRP.Range("H16").Value = 6750.154
$ws_CRP.Range("I16").Value = 5200.5
$ws_CRP.Range("J16").Value = 7438.8887
$ws_CRP.Range("K16").Value = 5200.5
$ws_CRP.Range("L16").Value = 7438.8887
$ws_CRP.Range("M16").Value = -4913.5
$ws_CRP.Range("N16").Value = -8012.8887
$ws_CRP.Range("H31").Value = 6398.65
$ws_CRP.Range("I31").Value = 3065.2368
$ws_CRP.Range("J31").Value = 12156.363
$ws_CRP.Range("K31").Value = 3065.2368
$ws_CRP.Range("L31").Value = 12156.363
$ws_CRP.Range("M31").Value = -2770.2368
$ws_CRP.Range("N31").Value = -12746.363
$ws_CRP.Range("H34").Value = 6398.65
$ws_CRP.Range("I34").Value = 3065.2368
$ws_CRP.Range("J34").Value = 12156.363
$ws_CRP.Range("K34").Value = 3065.2368
$ws_CRP.Range("L34").Value = 12156.363
$ws_CRP.Range("M34").Value = -2863.2368
$ws_CRP.Range("N34").Value = -12560.363
$ws_CRP.Range("H99").Value = 5713.8667
$ws_CRP.Range("I99").Value = 2836.75
$ws_CRP.Range("K99").Value = 2836.75
$ws_CRP.Range("M99").Value = -1338.75
$ws_CRP.Range("H113").Value = 6750.154
$ws_CRP.Range("I113").Value = 5200.5
$ws_CRP.Range("J113").Value = 7438.8887
$ws_CRP.Range("K113").Value = 5200.5
$ws_CRP.Range("L113").Value = 7438.8887
$ws_CRP.Range("M113").Value = -3030.5
$ws_CRP.Range("N113").Value = -11778.8887
$ws_CRP.Range("H126").Value = 5713.8667
$ws_CRP.Range("I126").Value = 2836.75
$ws_CRP.Range("K126").Value = 8510.25
$ws_CRP.Range("M126").Value = -6040.25
$ws_CRP.Range("H132").Value = 5790.646
$ws_CRP.Range("I132").Value = 3145.7058
$ws_CRP.Range("J132").Value = 12214.071
$ws_CRP.Range("K132").Value = 9437.117400000001
$ws_CRP.Range("L132").Value = 36642.213
$ws_CRP.Range("M132").Value = -6907.117400000001
$ws_CRP.Range("N132").Value = -41702.213
$ws_CRP.Range("H134").Value = 4169.0376
$ws_CRP.Range("I134").Value = 1816.1459
$ws_CRP.Range("K134").Value = 5448.4377
$ws_CRP.Range("M134").Value = -2913.4377
$ws_CUL.Range("H117").Value = 905.1429000000001
$ws_CUL.Range("J117").Value = 1405.3334
$ws_CUL.Range("L117").Value = 4216.0002
$ws_CUL.Range("N117").Value = -11100.0002
$ws_CUL.Range("H125").Value = 62505376
$ws_CUL.Range("I125").Value = 166669330
$ws_CUL.Range("J125").Value = 7000
$ws_CUL.Range("K125").Value = 500007990
$ws_CUL.Range("L125").Value = 21000
$ws_CUL.Range("M125").Value = -500003070
$ws_CUL.Range("N125").Value = -30840
$ws_CUL.Range("H129").Value = 937.5
$ws_CUL.Range("J129").Value = 2044
$ws_CUL.Range("L129").Value = 6132
$ws_CUL.Range("N129").Value = -16132
$ws_GSM.Range("H113").Value = 377759.22
$ws_GSM.Range("J113").Value = 9052.579
$ws_GSM.Range("L113").Value = 9052.579
$ws_GSM.Range("N113").Value = -13392.579
$ws_GSM.Range("H132").Value = 12658.941
$ws_GSM.Range("I132").Value = 6315.143
$ws_GSM.Range("J132").Value = 17099.6
$ws_GSM.Range("K132").Value = 18945.429
$ws_GSM.Range("L132").Value = 51298.8
$ws_GSM.Range("M132").Value = -16415.429
$ws_GSM.Range("N132").Value = -56358.8
$ws_LTW.Range("H55").Value = 398.08334
$ws_LTW.Range("I55").Value = 132.375
$ws_LTW.Range("J55").Value = 530.9375
$ws_LTW.Range("K55").Value = 132.375
$ws_LTW.Range("L55").Value = 530.9375
$ws_LTW.Range("M55").Value = 40.625
$ws_LTW.Range("N55").Value = -876.9375
$ws_LTW.Range("H61").Value = 5514.316
$ws_LTW.Range("I61").Value = 1721.7778
$ws_LTW.Range("K61").Value = 1721.7778
$ws_LTW.Range("M61").Value = -1519.7778
$ws_LTW.Range("H93").Value = 10428.571
$ws_LTW.Range("I93").Value = 9250
$ws_LTW.Range("K93").Value = 9250
$ws_LTW.Range("M93").Value = -8002
$ws_LTW.Range("H100").Value = 3750.2856
$ws_LTW.Range("J100").Value = 4597.3335
$ws_LTW.Range("L100").Value = 4597.3335
$ws_LTW.Range("N100").Value = -5679.3335
$ws_LTW.Range("H113").Value = 5514.316
$ws_LTW.Range("I113").Value = 1721.7778
$ws_LTW.Range("K113").Value = 1721.7778
$ws_LTW.Range("M113").Value = 448.2221999999999
$ws_LTW.Range("H122").Value = 6436.364
$ws_LTW.Range("I122").Value = 4759.7
$ws_LTW.Range("J122").Value = 7833.5835
$ws_LTW.Range("K122").Value = 14279.1
$ws_LTW.Range("L122").Value = 23500.7505
$ws_LTW.Range("M122").Value = -11829.1
$ws_LTW.Range("N122").Value = -28400.7505
$ws_LTW.Range("H130").Value = 59519
$ws_LTW.Range("J130").Value = 59519
$ws_LTW.Range("L130").Value = 59519
$ws_LTW.Range("N130").Value = -69559
$ws_WVR.Range("H15").Value = 32499
$ws_WVR.Range("H107").Value = 19608754
$ws_WVR.Range("I107").Value = 432.2
$ws_WVR.Range("K107").Value = 1296.6
$ws_WVR.Range("M107").Value = 623.4000000000001
$ws_WVR.Range("H132").Value = 12512169
$ws_WVR.Range("I132").Value = 20006268
$ws_WVR.Range("K132").Value = 60018804
$ws_WVR.Range("M132").Value = -60016274
